# Auto-generated: apply updated market price figures to each profession sheet's price/profit columns.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2273.8948
$ws.Range("I15").Value = 2273.8948
$ws.Range("K15").Value = 6821.6844
$ws.Range("M15").Value = -6652.6844
$ws.Range("H28").Value = 2125.375
$ws.Range("J28").Value = 3000.8
$ws.Range("L28").Value = 3000.8
$ws.Range("N28").Value = -3970.8
$ws.Range("H69").Value = 8176.091
$ws.Range("I69").Value = 6664.5
$ws.Range("K69").Value = 19993.5
$ws.Range("M69").Value = -19119.5
$ws.Range("H72").Value = 8176.091
$ws.Range("I72").Value = 6664.5
$ws.Range("K72").Value = 59980.5
$ws.Range("M72").Value = -55612.5
$ws.Range("H92").Value = 72985.78999999999
$ws.Range("I92").Value = 101446
$ws.Range("K92").Value = 101446
$ws.Range("M92").Value = -100198
$ws.Range("H100").Value = 3898.875
$ws.Range("I100").Value = 2858.8
$ws.Range("J100").Value = 5632.3335
$ws.Range("K100").Value = 2858.8
$ws.Range("L100").Value = 5632.3335
$ws.Range("M100").Value = -2317.8
$ws.Range("N100").Value = -6714.3335
$ws.Range("H106").Value = 7340.143
$ws.Range("I106").Value = 4396.75
$ws.Range("K106").Value = 4396.75
$ws.Range("M106").Value = -3765.75
$ws.Range("H107").Value = 1403
$ws.Range("I107").Value = 1135.8334
$ws.Range("K107").Value = 1135.8334
$ws.Range("M107").Value = 784.1666
$ws.Range("H116").Value = 9602.134
$ws.Range("I116").Value = 8536.223
$ws.Range("J116").Value = 11201
$ws.Range("K116").Value = 8536.223
$ws.Range("L116").Value = 11201
$ws.Range("M116").Value = -5094.223
$ws.Range("N116").Value = -18085
$ws.Range("H118").Value = 4277.8
$ws.Range("I118").Value = 463
$ws.Range("J118").Value = 10000
$ws.Range("K118").Value = 1389
$ws.Range("L118").Value = 30000
$ws.Range("M118").Value = 268
$ws.Range("N118").Value = -33314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 37385.234
$ws.Range("I32").Value = 38983.594
$ws.Range("K32").Value = 38983.594
$ws.Range("M32").Value = -38696.594
$ws.Range("H61").Value = 2918.7856
$ws.Range("I61").Value = 2758.6924
$ws.Range("K61").Value = 2758.6924
$ws.Range("M61").Value = -2546.6924
$ws.Range("H110").Value = 2299.4
$ws.Range("I110").Value = 2043.3077
$ws.Range("J110").Value = 2775
$ws.Range("K110").Value = 2043.3077
$ws.Range("L110").Value = 2775
$ws.Range("M110").Value = 1.692299999999932
$ws.Range("N110").Value = -6865
$ws.Range("H136").Value = 2918.7856
$ws.Range("I136").Value = 2758.6924
$ws.Range("K136").Value = 8276.0772
$ws.Range("M136").Value = -5726.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 49045.59
$ws.Range("I99").Value = 55213.316
$ws.Range("J99").Value = 9983.333000000001
$ws.Range("K99").Value = 55213.316
$ws.Range("L99").Value = 9983.333000000001
$ws.Range("M99").Value = -53715.316
$ws.Range("N99").Value = -12979.333
$ws.Range("H101").Value = 71841.5
$ws.Range("I101").Value = 73684
$ws.Range("J101").Value = 69999
$ws.Range("K101").Value = 73684
$ws.Range("L101").Value = 69999
$ws.Range("M101").Value = -70439
$ws.Range("N101").Value = -76489
$ws.Range("H102").Value = 13166.667
$ws.Range("I102").Value = 13166.667
$ws.Range("K102").Value = 13166.667
$ws.Range("M102").Value = -9921.666999999999
$ws.Range("H103").Value = 40910.6
$ws.Range("J103").Value = 40910.6
$ws.Range("L103").Value = 40910.6
$ws.Range("N103").Value = -43254.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2449
$ws.Range("J16").Value = 3248.5
$ws.Range("L16").Value = 3248.5
$ws.Range("N16").Value = -3822.5
$ws.Range("H58").Value = 47785.09
$ws.Range("I58").Value = 52202.9
$ws.Range("K58").Value = 52202.9
$ws.Range("M58").Value = -51999.9
$ws.Range("H86").Value = 6831.6665
$ws.Range("I86").Value = 6753.5
$ws.Range("K86").Value = 6753.5
$ws.Range("M86").Value = -5630.5
$ws.Range("H89").Value = 6831.6665
$ws.Range("I89").Value = 6753.5
$ws.Range("K89").Value = 33767.5
$ws.Range("M89").Value = -28151.5
$ws.Range("H99").Value = 4646.0713
$ws.Range("J99").Value = 4799.4
$ws.Range("L99").Value = 4799.4
$ws.Range("N99").Value = -7795.4
$ws.Range("H105").Value = 1802.3334
$ws.Range("I105").Value = 1802.3334
$ws.Range("K105").Value = 1802.3334
$ws.Range("M105").Value = -55.33339999999998
$ws.Range("H113").Value = 2449
$ws.Range("J113").Value = 3248.5
$ws.Range("L113").Value = 3248.5
$ws.Range("N113").Value = -7588.5
$ws.Range("H126").Value = 4646.0713
$ws.Range("J126").Value = 4799.4
$ws.Range("L126").Value = 14398.2
$ws.Range("N126").Value = -19338.2
$ws.Range("H132").Value = 1562.1538
$ws.Range("I132").Value = 1413.75
$ws.Range("J132").Value = 1799.6
$ws.Range("K132").Value = 4241.25
$ws.Range("L132").Value = 5398.799999999999
$ws.Range("M132").Value = -1711.25
$ws.Range("N132").Value = -10458.8
$ws.Range("H136").Value = 47785.09
$ws.Range("I136").Value = 52202.9
$ws.Range("K136").Value = 156608.7
$ws.Range("M136").Value = -154058.7
$ws.Range("H141").Value = 94999
$ws.Range("J141").Value = 94999
$ws.Range("L141").Value = 94999
$ws.Range("N141").Value = -105359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 226.72728
$ws.Range("J17").Value = 174.5
$ws.Range("L17").Value = 523.5
$ws.Range("N17").Value = -861.5
$ws.Range("H39").Value = 2232.1667
$ws.Range("I39").Value = 698.5
$ws.Range("K39").Value = 2095.5
$ws.Range("M39").Value = -1801.5
$ws.Range("H60").Value = 405.93332
$ws.Range("I60").Value = 349.125
$ws.Range("K60").Value = 1047.375
$ws.Range("M60").Value = -796.375
$ws.Range("H131").Value = 9654.521000000001
$ws.Range("J131").Value = 6698.5884
$ws.Range("L131").Value = 20095.7652
$ws.Range("N131").Value = -30175.7652

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1825.6
$ws.Range("I7").Value = 1825.6
$ws.Range("K7").Value = 1825.6
$ws.Range("M7").Value = -1713.6
$ws.Range("H40").Value = 6497
$ws.Range("I40").Value = 2400.6
$ws.Range("K40").Value = 2400.6
$ws.Range("M40").Value = -2264.6
$ws.Range("H93").Value = 2532.3333
$ws.Range("I93").Value = 1986.375
$ws.Range("K93").Value = 1986.375
$ws.Range("M93").Value = -738.375
$ws.Range("H100").Value = 4470.48
$ws.Range("I100").Value = 4224.1177
$ws.Range("J100").Value = 4994
$ws.Range("K100").Value = 4224.1177
$ws.Range("L100").Value = 4994
$ws.Range("M100").Value = -3683.1177
$ws.Range("N100").Value = -6076
$ws.Range("H126").Value = 1825.6
$ws.Range("I126").Value = 1825.6
$ws.Range("K126").Value = 5476.799999999999
$ws.Range("M126").Value = -3006.799999999999
$ws.Range("H129").Value = 79999
$ws.Range("J129").Value = 79999
$ws.Range("L129").Value = 79999
$ws.Range("N129").Value = -89999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 831.46155
$ws.Range("I113").Value = 375.75
$ws.Range("K113").Value = 1127.25
$ws.Range("M113").Value = 1042.75
$ws.Range("H126").Value = 59972.445
$ws.Range("I126").Value = 75549.92999999999
$ws.Range("K126").Value = 226649.79
$ws.Range("M126").Value = -224179.79
